$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.5429463333333333
$ws.Range("H2").Value = 1.628839
$ws.Range("I2").Value = 0.04659251079363984
$ws.Range("J2").Value = 0.04659251079363985
$ws.Range("M2").Value = 68.63737500000001
$ws.Range("N2").Value = 205.912125
$ws.Range("O2").Value = 0.5415701538216162
$ws.Range("P2").Value = 0.5415701538216162
$ws.Range("Q2").Value = 37.266411085875
$ws.Range("R2").Value = 335.397699772875
$ws.Range("S2").Value = 0.02523311323744684
$ws.Range("T2").Value = 0.02523311323744685
$ws.Range("G3").Value = 0.5429463333333333
$ws.Range("H3").Value = 1.628839
$ws.Range("I3").Value = 0.04659251079363984
$ws.Range("J3").Value = 0.04659251079363985
$ws.Range("O3").Value = 0.08718851262838957
$ws.Range("P3").Value = 0.08718851262838957
$ws.Range("Q3").Value = 5.999597523326222
$ws.Range("R3").Value = 53.996377709936
$ws.Range("S3").Value = 0.004062331715719645
$ws.Range("T3").Value = 0.004062331715719646
$ws.Range("G4").Value = 0.5429463333333333
$ws.Range("H4").Value = 1.628839
$ws.Range("I4").Value = 0.04659251079363984
$ws.Range("J4").Value = 0.04659251079363985
$ws.Range("M4").Value = 16.21089566666667
$ws.Range("N4").Value = 48.632687
$ws.Range("O4").Value = 0.1279089892319285
$ws.Range("P4").Value = 0.1279089892319285
$ws.Range("Q4").Value = 8.80164636226589
$ws.Range("R4").Value = 79.214817260393
$ws.Range("S4").Value = 0.005959600961392191
$ws.Range("T4").Value = 0.005959600961392192
$ws.Range("G5").Value = 0.5429463333333333
$ws.Range("H5").Value = 1.628839
$ws.Range("I5").Value = 0.04659251079363984
$ws.Range("J5").Value = 0.04659251079363985
$ws.Range("M5").Value = 20.32546233333333
$ws.Range("N5").Value = 60.976387
$ws.Range("O5").Value = 0.1603741949973873
$ws.Range("P5").Value = 0.1603741949973873
$ws.Range("Q5").Value = 11.03563524718811
$ws.Range("R5").Value = 99.320717224693
$ws.Range("S5").Value = 0.00747223641143707
$ws.Range("T5").Value = 0.007472236411437072
$ws.Range("G6").Value = 0.5429463333333333
$ws.Range("H6").Value = 1.628839
$ws.Range("I6").Value = 0.04659251079363984
$ws.Range("J6").Value = 0.04659251079363985
$ws.Range("M6").Value = 10.513928
$ws.Range("N6").Value = 31.541784
$ws.Range("O6").Value = 0.08295814932067838
$ws.Range("P6").Value = 0.08295814932067838
$ws.Range("Q6").Value = 5.708498656530667
$ws.Range("R6").Value = 51.376487908776
$ws.Range("S6").Value = 0.003865228467644094
$ws.Range("T6").Value = 0.003865228467644094
$ws.Range("I7").Value = 0.8858267105024722
$ws.Range("J7").Value = 0.8858267105024723
$ws.Range("M7").Value = 68.63737500000001
$ws.Range("N7").Value = 205.912125
$ws.Range("O7").Value = 0.5415701538216162
$ws.Range("P7").Value = 0.5415701538216162
$ws.Range("Q7").Value = 708.516922186125
$ws.Range("R7").Value = 6376.652299675125
$ws.Range("S7").Value = 0.4797373078661202
$ws.Range("T7").Value = 0.4797373078661202
$ws.Range("I8").Value = 0.8858267105024722
$ws.Range("J8").Value = 0.8858267105024723
$ws.Range("O8").Value = 0.08718851262838957
$ws.Range("P8").Value = 0.08718851262838957
$ws.Range("S8").Value = 0.07723391333520958
$ws.Range("T8").Value = 0.0772339133352096
$ws.Range("I9").Value = 0.8858267105024722
$ws.Range("J9").Value = 0.8858267105024723
$ws.Range("M9").Value = 16.21089566666667
$ws.Range("N9").Value = 48.632687
$ws.Range("O9").Value = 0.1279089892319285
$ws.Range("P9").Value = 0.1279089892319285
$ws.Range("Q9").Value = 167.3387699285857
$ws.Range("R9").Value = 1506.048929357271
$ws.Range("S9").Value = 0.1133051991750154
$ws.Range("T9").Value = 0.1133051991750154
$ws.Range("I10").Value = 0.8858267105024722
$ws.Range("J10").Value = 0.8858267105024723
$ws.Range("M10").Value = 20.32546233333333
$ws.Range("N10").Value = 60.976387
$ws.Range("O10").Value = 0.1603741949973873
$ws.Range("P10").Value = 0.1603741949973873
$ws.Range("Q10").Value = 209.8118410621524
$ws.Range("R10").Value = 1888.306569559371
$ws.Range("S10").Value = 0.1420637456040177
$ws.Range("T10").Value = 0.1420637456040177
$ws.Range("I11").Value = 0.8858267105024722
$ws.Range("J11").Value = 0.8858267105024723
$ws.Range("M11").Value = 10.513928
$ws.Range("N11").Value = 31.541784
$ws.Range("O11").Value = 0.08295814932067838
$ws.Range("P11").Value = 0.08295814932067838
$ws.Range("Q11").Value = 108.531188826008
$ws.Range("R11").Value = 976.7806994340719
$ws.Range("S11").Value = 0.07348654452210943
$ws.Range("T11").Value = 0.07348654452210944
$ws.Range("G12").Value = 0.7875243333333334
$ws.Range("H12").Value = 2.362573
$ws.Range("I12").Value = 0.06758077870388791
$ws.Range("J12").Value = 0.06758077870388793
$ws.Range("M12").Value = 68.63737500000001
$ws.Range("N12").Value = 205.912125
$ws.Range("O12").Value = 0.5415701538216162
$ws.Range("P12").Value = 0.5415701538216162
$ws.Range("Q12").Value = 54.05360298862501
$ws.Range("R12").Value = 486.482426897625
$ws.Range("S12").Value = 0.03659973271804919
$ws.Range("T12").Value = 0.03659973271804919
$ws.Range("G13").Value = 0.7875243333333334
$ws.Range("H13").Value = 2.362573
$ws.Range("I13").Value = 0.06758077870388791
$ws.Range("J13").Value = 0.06758077870388793
$ws.Range("O13").Value = 0.08718851262838957
$ws.Range("P13").Value = 0.08718851262838957
$ws.Range("Q13").Value = 8.702202685150223
$ws.Range("R13").Value = 78.31982416635201
$ws.Range("S13").Value = 0.005892267577460333
$ws.Range("T13").Value = 0.005892267577460334
$ws.Range("G14").Value = 0.7875243333333334
$ws.Range("H14").Value = 2.362573
$ws.Range("I14").Value = 0.06758077870388791
$ws.Range("J14").Value = 0.06758077870388793
$ws.Range("M14").Value = 16.21089566666667
$ws.Range("N14").Value = 48.632687
$ws.Range("O14").Value = 0.1279089892319285
$ws.Range("P14").Value = 0.1279089892319285
$ws.Range("Q14").Value = 12.76647480262789
$ws.Range("R14").Value = 114.898273223651
$ws.Range("S14").Value = 0.008644189095520942
$ws.Range("T14").Value = 0.008644189095520944
$ws.Range("G15").Value = 0.7875243333333334
$ws.Range("H15").Value = 2.362573
$ws.Range("I15").Value = 0.06758077870388791
$ws.Range("J15").Value = 0.06758077870388793
$ws.Range("M15").Value = 20.32546233333333
$ws.Range("N15").Value = 60.976387
$ws.Range("O15").Value = 0.1603741949973873
$ws.Range("P15").Value = 0.1603741949973873
$ws.Range("Q15").Value = 16.00679617375011
$ws.Range("R15").Value = 144.061165563751
$ws.Range("S15").Value = 0.0108382129819326
$ws.Range("T15").Value = 0.0108382129819326
$ws.Range("G16").Value = 0.7875243333333334
$ws.Range("H16").Value = 2.362573
$ws.Range("I16").Value = 0.06758077870388791
$ws.Range("J16").Value = 0.06758077870388793
$ws.Range("M16").Value = 10.513928
$ws.Range("N16").Value = 31.541784
$ws.Range("O16").Value = 0.08295814932067838
$ws.Range("P16").Value = 0.08295814932067838
$ws.Range("Q16").Value = 8.279974138914667
$ws.Range("R16").Value = 74.51976725023201
$ws.Range("S16").Value = 0.005606376330924855
$ws.Range("T16").Value = 0.005606376330924857
